$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of workhours data following Unity's official WheelCollider tutorial
$ws.Range("B3").Value = 0.30902777777777779
$ws.Range("C3").Value = 0.31597222222222221
$ws.Range("D3").Value = "Unity WheelCollider documentation"

$ws.Range("B4").Value = 0.31597222222222221
$ws.Range("C4").Value = 0.33333333333333331
$ws.Range("D4").Value = "Unity WheelCollider tutorial"

$ws.Range("B5").Value = 0.34027777777777773
$ws.Range("C5").Value = 0.375
$ws.Range("D5").Value = "Unity WheelCollider tutorial"

# Apply time number format (h:mm, builtin numFmtId 20) to the new From/To cells
$ws.Range("B3:C5").NumberFormat = "h:mm"

# Match the saved selection state
$ws.Range("G5").Select()
